# Append three new payment-log rows (31-33) to the "Payments" sheet, matching
# the rows produced by the site's "success page GET" handler being invoked
# multiple times for the same checkout flow.
#
# All columns except O (Amount) hold values that must be preserved verbatim as
# TEXT, even when they look numeric (dates typed as "1/1/2025", zero-padded
# transaction ids, batch codes like "21", phone numbers, etc.) - exactly like
# the existing rows already on the sheet. Excel would otherwise auto-convert
# those into numbers/dates on entry, so each text cell is briefly switched to
# the "@" (Text) number format before the value is typed in, then the
# temporary formatting is cleared again so the cell's style matches the rest
# of the sheet (plain/default style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-PaymentRow {
    param($r, $vals)

    # Columns that must stay literal text (everything except N=FeeSem, which
    # is left blank, and O=Amount, which is a real number).
    $textCols = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 16, 17)

    foreach ($c in $textCols) {
        $ws.Cells.Item($r, $c).NumberFormat = "@"
    }

    for ($i = 0; $i -lt $vals.Count; $i++) {
        $c = $i + 1
        if ($c -eq 14) {
            # N / FeeSem - left blank for these rows, same as the diff's <v/>
            continue
        }
        $ws.Cells.Item($r, $c).Value = $vals[$i]
    }

    foreach ($c in $textCols) {
        $ws.Cells.Item($r, $c).ClearFormats()
    }
}

# Columns: A Date, B Time, C TransactionId, D Razorpay_Order_Id, E Batch,
#          F Roll, G Name, H Branch, I Section, J Phone, K Email, L FeeType,
#          M FeeYear, N FeeSem, O Amount, P Method, Q Status

$row31 = @(
    "1/1/2025", "4:57:58 pm", "010125165758", "order_PeA8Jf1J1MAb0G",
    "21", "21B81A05V9", "SAMRATH REDDY", "CSE", "E", "+917981455290",
    "samrathreddy04@gmail.com", "CollegeFee", "IV", "",
    120000, "wallet", "Verification in progress..."
)

$row32 = @(
    "1/1/2025", "5:27:06 pm", "010125172706", "order_PeAd37PCqw6j8f",
    "21", "21B81A05V9", "SAMRATH REDDY", "CSE", "E", "+917981455290",
    "samrathreddy04@gmail.com", "CollegeFee", "IV", "",
    120000, "wallet", "Verification in progress..."
)

$row33 = @(
    "1/1/2025", "5:29:27 pm", "010125172927", "order_PeAfX0B2PHevt5",
    "21", "21B81A05V9", "SAMRATH REDDY", "CSE", "E", "+917981455290",
    "samrathreddy04@gmail.com", "CollegeFee", "I", "",
    120000, "wallet", "Rejected"
)

Add-PaymentRow 31 $row31
Add-PaymentRow 32 $row32
Add-PaymentRow 33 $row33
